$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2022-08-12 20:58:15"

for ($row = 2; $row -le 73; $row++) {
    $ws.Cells.Item($row, 15).Value = $newTimestamp
}
